$d = $word.ActiveDocument

$d.Content.Find.Execute("56-25=31", $true, $true, $false, $false, $false, $true, 1, $false, "1+97=98", 2) | Out-Null
$d.Content.Find.Execute("72-66=6", $true, $true, $false, $false, $false, $true, 1, $false, "75-67=8", 2) | Out-Null
$d.Content.Find.Execute("33+5=38", $true, $true, $false, $false, $false, $true, 1, $false, "29+7=36", 2) | Out-Null
$d.Content.Find.Execute("5-2=3", $true, $true, $false, $false, $false, $true, 1, $false, "62-41=21", 2) | Out-Null
$d.Content.Find.Execute("86+11=97", $true, $true, $false, $false, $false, $true, 1, $false, "35+60=95", 2) | Out-Null
$d.Content.Find.Execute("74-1=73", $true, $true, $false, $false, $false, $true, 1, $false, "80-73=7", 2) | Out-Null
$d.Content.Find.Execute("8+38=46", $true, $true, $false, $false, $false, $true, 1, $false, "57-8=49", 2) | Out-Null
$d.Content.Find.Execute("94-28=66", $true, $true, $false, $false, $false, $true, 1, $false, "82+10=92", 2) | Out-Null
$d.Content.Find.Execute("45+3=48", $true, $true, $false, $false, $false, $true, 1, $false, "49+1=50", 2) | Out-Null
$d.Content.Find.Execute("59+17=76", $true, $true, $false, $false, $false, $true, 1, $false, "31+45=76", 2) | Out-Null
$d.Content.Find.Execute("87-67=20", $true, $true, $false, $false, $false, $true, 1, $false, "20+27=47", 2) | Out-Null
$d.Content.Find.Execute("8+70=78", $true, $true, $false, $false, $false, $true, 1, $false, "85-53=32", 2) | Out-Null
$d.Content.Find.Execute("78+9=87", $true, $true, $false, $false, $false, $true, 1, $false, "48-42=6", 2) | Out-Null
$d.Content.Find.Execute("22+74=96", $true, $true, $false, $false, $false, $true, 1, $false, "66-9=57", 2) | Out-Null
$d.Content.Find.Execute("8+20=28", $true, $true, $false, $false, $false, $true, 1, $false, "22+12=34", 2) | Out-Null
$d.Content.Find.Execute("13+45=58", $true, $true, $false, $false, $false, $true, 1, $false, "85-43=42", 2) | Out-Null
$d.Content.Find.Execute("65+7=72", $true, $true, $false, $false, $false, $true, 1, $false, "85-75=10", 2) | Out-Null
$d.Content.Find.Execute("85-3=82", $true, $true, $false, $false, $false, $true, 1, $false, "53-17=36", 2) | Out-Null
$d.Content.Find.Execute("34+43=77", $true, $true, $false, $false, $false, $true, 1, $false, "6+2=8", 2) | Out-Null
$d.Content.Find.Execute("68-29=39", $true, $true, $false, $false, $false, $true, 1, $false, "38+47=85", 2) | Out-Null
$d.Content.Find.Execute("8+10=18", $true, $true, $false, $false, $false, $true, 1, $false, "44-4=40", 2) | Out-Null
$d.Content.Find.Execute("33-26=7", $true, $true, $false, $false, $false, $true, 1, $false, "55-51=4", 2) | Out-Null
$d.Content.Find.Execute("68-23=45", $true, $true, $false, $false, $false, $true, 1, $false, "27-20=7", 2) | Out-Null
$d.Content.Find.Execute("38+37=75", $true, $true, $false, $false, $false, $true, 1, $false, "93-60=33", 2) | Out-Null
$d.Content.Find.Execute("62+23=85", $true, $true, $false, $false, $false, $true, 1, $false, "86-48=38", 2) | Out-Null
$d.Content.Find.Execute("12+3=15", $true, $true, $false, $false, $false, $true, 1, $false, "41-2=39", 2) | Out-Null
$d.Content.Find.Execute("31+58=89", $true, $true, $false, $false, $false, $true, 1, $false, "51+48=99", 2) | Out-Null
$d.Content.Find.Execute("99-77=22", $true, $true, $false, $false, $false, $true, 1, $false, "2+3=5", 2) | Out-Null
$d.Content.Find.Execute("73-33=40", $true, $true, $false, $false, $false, $true, 1, $false, "82+1=83", 2) | Out-Null
$d.Content.Find.Execute("99-65=34", $true, $true, $false, $false, $false, $true, 1, $false, "84-10=74", 2) | Out-Null
$d.Content.Find.Execute("92-31=61", $true, $true, $false, $false, $false, $true, 1, $false, "53+39=92", 2) | Out-Null
$d.Content.Find.Execute("55-16=39", $true, $true, $false, $false, $false, $true, 1, $false, "79-74=5", 2) | Out-Null
$d.Content.Find.Execute("21-15=6", $true, $true, $false, $false, $false, $true, 1, $false, "32-23=9", 2) | Out-Null
$d.Content.Find.Execute("6+41=47", $true, $true, $false, $false, $false, $true, 1, $false, "83-24=59", 2) | Out-Null
$d.Content.Find.Execute("68-45=23", $true, $true, $false, $false, $false, $true, 1, $false, "25-7=18", 2) | Out-Null
$d.Content.Find.Execute("42-17=25", $true, $true, $false, $false, $false, $true, 1, $false, "13+78=91", 2) | Out-Null
$d.Content.Find.Execute("95-61=34", $true, $true, $false, $false, $false, $true, 1, $false, "48-20=28", 2) | Out-Null
$d.Content.Find.Execute("41+47=88", $true, $true, $false, $false, $false, $true, 1, $false, "24+50=74", 2) | Out-Null
$d.Content.Find.Execute("72+16=88", $true, $true, $false, $false, $false, $true, 1, $false, "66-46=20", 2) | Out-Null
$d.Content.Find.Execute("15+46=61", $true, $true, $false, $false, $false, $true, 1, $false, "56+37=93", 2) | Out-Null
$d.Content.Find.Execute("2+57=59", $true, $true, $false, $false, $false, $true, 1, $false, "0+27=27", 2) | Out-Null
$d.Content.Find.Execute("70+24=94", $true, $true, $false, $false, $false, $true, 1, $false, "92-22=70", 2) | Out-Null
$d.Content.Find.Execute("54-1=53", $true, $true, $false, $false, $false, $true, 1, $false, "41+39=80", 2) | Out-Null
$d.Content.Find.Execute("17-16=1", $true, $true, $false, $false, $false, $true, 1, $false, "42-29=13", 2) | Out-Null
$d.Content.Find.Execute("30+26=56", $true, $true, $false, $false, $false, $true, 1, $false, "84-22=62", 2) | Out-Null
$d.Content.Find.Execute("29-12=17", $true, $true, $false, $false, $false, $true, 1, $false, "22+18=40", 2) | Out-Null
$d.Content.Find.Execute("92-48=44", $true, $true, $false, $false, $false, $true, 1, $false, "9+54=63", 2) | Out-Null
$d.Content.Find.Execute("24-8=16", $true, $true, $false, $false, $false, $true, 1, $false, "45-21=24", 2) | Out-Null
$d.Content.Find.Execute("78-77=1", $true, $true, $false, $false, $false, $true, 1, $false, "68+31=99", 2) | Out-Null
$d.Content.Find.Execute("13+21=34", $true, $true, $false, $false, $false, $true, 1, $false, "67-67=0", 2) | Out-Null
$d.Content.Find.Execute("41-38=3", $true, $true, $false, $false, $false, $true, 1, $false, "6+21=27", 2) | Out-Null
$d.Content.Find.Execute("95-54=41", $true, $true, $false, $false, $false, $true, 1, $false, "58-11=47", 2) | Out-Null
$d.Content.Find.Execute("25-11=14", $true, $true, $false, $false, $false, $true, 1, $false, "32-31=1", 2) | Out-Null
$d.Content.Find.Execute("44-19=25", $true, $true, $false, $false, $false, $true, 1, $false, "61-24=37", 2) | Out-Null
$d.Content.Find.Execute("48+18=66", $true, $true, $false, $false, $false, $true, 1, $false, "8+32=40", 2) | Out-Null
$d.Content.Find.Execute("23+23=46", $true, $true, $false, $false, $false, $true, 1, $false, "14+65=79", 2) | Out-Null
$d.Content.Find.Execute("73-19=54", $true, $true, $false, $false, $false, $true, 1, $false, "66-48=18", 2) | Out-Null
$d.Content.Find.Execute("77-16=61", $true, $true, $false, $false, $false, $true, 1, $false, "48-11=37", 2) | Out-Null
$d.Content.Find.Execute("3+50=53", $true, $true, $false, $false, $false, $true, 1, $false, "82-32=50", 2) | Out-Null
$d.Content.Find.Execute("72+1=73", $true, $true, $false, $false, $false, $true, 1, $false, "41-11=30", 2) | Out-Null
$d.Content.Find.Execute("20+38=58", $true, $true, $false, $false, $false, $true, 1, $false, "8+48=56", 2) | Out-Null
$d.Content.Find.Execute("26+53=79", $true, $true, $false, $false, $false, $true, 1, $false, "94-31=63", 2) | Out-Null
$d.Content.Find.Execute("5+84=89", $true, $true, $false, $false, $false, $true, 1, $false, "24+53=77", 2) | Out-Null
$d.Content.Find.Execute("85-84=1", $true, $true, $false, $false, $false, $true, 1, $false, "26-12=14", 2) | Out-Null
$d.Content.Find.Execute("17+52=69", $true, $true, $false, $false, $false, $true, 1, $false, "93-61=32", 2) | Out-Null
$d.Content.Find.Execute("11+14=25", $true, $true, $false, $false, $false, $true, 1, $false, "12+66=78", 2) | Out-Null
$d.Content.Find.Execute("16+62=78", $true, $true, $false, $false, $false, $true, 1, $false, "80-67=13", 2) | Out-Null
$d.Content.Find.Execute("81+10=91", $true, $true, $false, $false, $false, $true, 1, $false, "6+10=16", 2) | Out-Null
$d.Content.Find.Execute("1+84=85", $true, $true, $false, $false, $false, $true, 1, $false, "72-25=47", 2) | Out-Null
$d.Content.Find.Execute("68-4=64", $true, $true, $false, $false, $false, $true, 1, $false, "20+21=41", 2) | Out-Null
$d.Content.Find.Execute("60-42=18", $true, $true, $false, $false, $false, $true, 1, $false, "12+56=68", 2) | Out-Null
$d.Content.Find.Execute("11+39=50", $true, $true, $false, $false, $false, $true, 1, $false, "61+0=61", 2) | Out-Null
$d.Content.Find.Execute("94-67=27", $true, $true, $false, $false, $false, $true, 1, $false, "39+5=44", 2) | Out-Null
$d.Content.Find.Execute("14+21=35", $true, $true, $false, $false, $false, $true, 1, $false, "48+36=84", 2) | Out-Null
$d.Content.Find.Execute("65-30=35", $true, $true, $false, $false, $false, $true, 1, $false, "34+29=63", 2) | Out-Null
$d.Content.Find.Execute("42-34=8", $true, $true, $false, $false, $false, $true, 1, $false, "2+8=10", 2) | Out-Null
$d.Content.Find.Execute("53+40=93", $true, $true, $false, $false, $false, $true, 1, $false, "18+78=96", 2) | Out-Null
$d.Content.Find.Execute("75-65=10", $true, $true, $false, $false, $false, $true, 1, $false, "34+2=36", 2) | Out-Null
$d.Content.Find.Execute("13+1=14", $true, $true, $false, $false, $false, $true, 1, $false, "5+42=47", 2) | Out-Null
$d.Content.Find.Execute("61-59=2", $true, $true, $false, $false, $false, $true, 1, $false, "47+19=66", 2) | Out-Null
$d.Content.Find.Execute("17-1=16", $true, $true, $false, $false, $false, $true, 1, $false, "1+30=31", 2) | Out-Null
$d.Content.Find.Execute("37+6=43", $true, $true, $false, $false, $false, $true, 1, $false, "53-20=33", 2) | Out-Null
$d.Content.Find.Execute("27+40=67", $true, $true, $false, $false, $false, $true, 1, $false, "7+90=97", 2) | Out-Null
$d.Content.Find.Execute("70+12=82", $true, $true, $false, $false, $false, $true, 1, $false, "70-17=53", 2) | Out-Null
$d.Content.Find.Execute("46-21=25", $true, $true, $false, $false, $false, $true, 1, $false, "76+15=91", 2) | Out-Null
$d.Content.Find.Execute("39-13=26", $true, $true, $false, $false, $false, $true, 1, $false, "94-1=93", 2) | Out-Null
$d.Content.Find.Execute("58+28=86", $true, $true, $false, $false, $false, $true, 1, $false, "30+33=63", 2) | Out-Null
$d.Content.Find.Execute("86-45=41", $true, $true, $false, $false, $false, $true, 1, $false, "73-2=71", 2) | Out-Null
$d.Content.Find.Execute("18+53=71", $true, $true, $false, $false, $false, $true, 1, $false, "33+43=76", 2) | Out-Null
$d.Content.Find.Execute("95+0=95", $true, $true, $false, $false, $false, $true, 1, $false, "9+1=10", 2) | Out-Null
$d.Content.Find.Execute("55+43=98", $true, $true, $false, $false, $false, $true, 1, $false, "46+42=88", 2) | Out-Null
$d.Content.Find.Execute("47+35=82", $true, $true, $false, $false, $false, $true, 1, $false, "99-35=64", 2) | Out-Null
$d.Content.Find.Execute("76-19=57", $true, $true, $false, $false, $false, $true, 1, $false, "16+83=99", 2) | Out-Null
$d.Content.Find.Execute("1+40=41", $true, $true, $false, $false, $false, $true, 1, $false, "31+23=54", 2) | Out-Null
$d.Content.Find.Execute("35-28=7", $true, $true, $false, $false, $false, $true, 1, $false, "4+88=92", 2) | Out-Null
$d.Content.Find.Execute("18+5=23", $true, $true, $false, $false, $false, $true, 1, $false, "90-17=73", 2) | Out-Null
$d.Content.Find.Execute("29+67=96", $true, $true, $false, $false, $false, $true, 1, $false, "50-23=27", 2) | Out-Null
$d.Content.Find.Execute("76-76=0", $true, $true, $false, $false, $false, $true, 1, $false, "13+66=79", 2) | Out-Null
$d.Content.Find.Execute("33+55=88", $true, $true, $false, $false, $false, $true, 1, $false, "24+15=39", 2) | Out-Null
$d.Content.Find.Execute("48-9=39", $true, $true, $false, $false, $false, $true, 1, $false, "3+21=24", 2) | Out-Null
